$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 311
$ws.Range("C4").Value = "abc"
$ws.Range("D4").Value = "tishya@gmail.com"
$ws.Range("E4").Value = "globalTiger"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "55667788"
$ws.Range("G4").Value = "L2_selected"
